$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: clear the y_0_forecast (C2) and y_1_forecast (E2) cells - no longer computable
$ws.Range("C2").ClearContents()
$ws.Range("E2").ClearContents()

# Row 3: clear y_0_forecast (C3); refine y_1_forecast (E3) precision
$ws.Range("C3").ClearContents()
$ws.Range("E3").Value = -1.563633964192079

# Row 4: refine y_0_forecast (C4) precision
$ws.Range("C4").Value = -4.774178217057756

# Row 6: refine y_1_forecast (E6) precision
$ws.Range("E6").Value = 2.114249845651872

# Row 7: refine y_0_forecast (C7) precision
$ws.Range("C7").Value = 1.239479831392831

# Row 8: refine y_0_forecast (C8) precision
$ws.Range("C8").Value = 0.2379616621360992

# Row 10: refine y_0_forecast (C10) and y_1_forecast (E10) precision
$ws.Range("C10").Value = 1.470039379455734
$ws.Range("E10").Value = 1.339087911421122

# Row 11: refine y_0_forecast (C11) and y_1_forecast (E11) precision
$ws.Range("C11").Value = 1.638797242243228
$ws.Range("E11").Value = 1.006353890555212

# Row 13: refine y_1_forecast (E13) precision
$ws.Range("E13").Value = -0.301339632123987

# Row 14: refine y_1_forecast (E14) precision
$ws.Range("E14").Value = 0.2691345740890139

# Row 15: refine y_1_forecast (E15) precision
$ws.Range("E15").Value = 23.52713729381606

# Row 16: refine y_0_forecast (C16) and y_1_forecast (E16) precision
$ws.Range("C16").Value = 1.099928004397577
$ws.Range("E16").Value = 6.182044950645027

# Row 17: refine y_0_forecast (C17) precision
$ws.Range("C17").Value = 2.310042359896247

# Row 18: refine y_1_forecast (E18) precision
$ws.Range("E18").Value = -0.3513551123189074

# Row 19: refine y_0_forecast (C19) and y_1_forecast (E19) precision
$ws.Range("C19").Value = -0.3101476031197037
$ws.Range("E19").Value = 0.2561130241983456
